$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col=2; Value=1.741367521909297},
    @{Row=2; Col=3; Value=0.2277835707313614},
    @{Row=2; Col=4; Value=0.5176247577602169},
    @{Row=2; Col=5; Value=0.174230358010373},
    @{Row=2; Col=7; Value=0.9400624967631472},
    @{Row=2; Col=8; Value=0.988284131320242},
    @{Row=2; Col=10; Value=0.07360246527502845},
    @{Row=2; Col=12; Value=0.4269343701549246},
    @{Row=2; Col=15; Value=3.898970018771109},
    @{Row=3; Col=2; Value=1.608189009190994},
    @{Row=3; Col=3; Value=0.2141702314516181},
    @{Row=3; Col=4; Value=0.5175166576501056},
    @{Row=3; Col=5; Value=0.1755950638993848},
    @{Row=3; Col=7; Value=0.948389781801076},
    @{Row=3; Col=8; Value=0.9975914093568718},
    @{Row=3; Col=10; Value=0.07408964589558664},
    @{Row=3; Col=12; Value=0.4165757608473797},
    @{Row=3; Col=15; Value=3.935756533715278},
    @{Row=4; Col=2; Value=1.526517907096661},
    @{Row=4; Col=3; Value=0.2057972729158735},
    @{Row=4; Col=4; Value=0.5176918674961115},
    @{Row=4; Col=5; Value=0.176506218026514},
    @{Row=4; Col=7; Value=0.9542189710910591},
    @{Row=4; Col=8; Value=1.003823009224973},
    @{Row=4; Col=10; Value=0.07440738117021528},
    @{Row=4; Col=12; Value=0.4103499208600709},
    @{Row=4; Col=15; Value=3.960931741061003},
    @{Row=5; Col=2; Value=1.493263847576486},
    @{Row=5; Col=3; Value=0.2023818992774409},
    @{Row=5; Col=4; Value=0.5178240902226179},
    @{Row=5; Col=5; Value=0.1768959484964387},
    @{Row=5; Col=7; Value=0.9567742450055334},
    @{Row=5; Col=8; Value=1.006492406109892},
    @{Row=5; Col=10; Value=0.07454154950090164},
    @{Row=5; Col=12; Value=0.4078468370327641},
    @{Row=5; Col=15; Value=3.97184106891747},
    @{Row=6; Col=2; Value=1.487743766240783},
    @{Row=6; Col=3; Value=0.2018145863943062},
    @{Row=6; Col=4; Value=0.5178497214986209},
    @{Row=6; Col=5; Value=0.1769617764123783},
    @{Row=6; Col=7; Value=0.957209399979881},
    @{Row=6; Col=8; Value=1.006943507036908},
    @{Row=6; Col=10; Value=0.07456411153845988},
    @{Row=6; Col=12; Value=0.4074332611992588},
    @{Row=6; Col=15; Value=3.973691809506377},
    @{Row=7; Col=2; Value=1.526069316669464},
    @{Row=7; Col=3; Value=0.2057512250891307},
    @{Row=7; Col=4; Value=0.5176934043261952},
    @{Row=7; Col=5; Value=0.1765113994330694},
    @{Row=7; Col=7; Value=0.9542527046273932},
    @{Row=7; Col=8; Value=1.003858483416749},
    @{Row=7; Col=10; Value=0.0744091716131754},
    @{Row=7; Col=12; Value=0.4103160254732074},
    @{Row=7; Col=15; Value=3.961076236111339},
    @{Row=8; Col=2; Value=1.695427852499733},
    @{Row=8; Col=3; Value=0.2230928212405558},
    @{Row=8; Col=4; Value=0.5175373878271756},
    @{Row=8; Col=5; Value=0.1746857216922333},
    @{Row=8; Col=7; Value=0.9427849381505098},
    @{Row=8; Col=8; Value=0.9913859949572696},
    @{Row=8; Col=10; Value=0.0737665900633937},
    @{Row=8; Col=12; Value=0.4233349471945473},
    @{Row=8; Col=15; Value=3.911116490138369},
    @{Row=9; Col=2; Value=2.02826327416949},
    @{Row=9; Col=3; Value=0.2569758062508072},
    @{Row=9; Col=4; Value=0.5191455953251563},
    @{Row=9; Col=5; Value=0.1716858302260196},
    @{Row=9; Col=7; Value=0.9259917476473163},
    @{Row=9; Col=8; Value=0.9710292684470545},
    @{Row=9; Col=10; Value=0.07265363752805243},
    @{Row=9; Col=12; Value=0.4499240255853607},
    @{Row=9; Col=15; Value=3.833708764362399},
    @{Row=10; Col=2; Value=2.273159735430852},
    @{Row=10; Col=3; Value=0.2817829649070234},
    @{Row=10; Col=4; Value=0.5214913145313318},
    @{Row=10; Col=5; Value=0.1698346148322685},
    @{Row=10; Col=7; Value=0.917143365231226},
    @{Row=10; Col=8; Value=0.9585745649043531},
    @{Row=10; Col=10; Value=0.0719250060335721},
    @{Row=10; Col=12; Value=0.4700972468479563},
    @{Row=10; Col=15; Value=3.789412973155748},
    @{Row=11; Col=2; Value=2.384633288481552},
    @{Row=11; Col=3; Value=0.2930474080404508},
    @{Row=11; Col=4; Value=0.5228107115820819},
    @{Row=11; Col=5; Value=0.1690688625677179},
    @{Row=11; Col=7; Value=0.9138795819521306},
    @{Row=11; Col=8; Value=0.9534519602173361},
    @{Row=11; Col=10; Value=0.07161273302147997},
    @{Row=11; Col=12; Value=0.4794117624662704},
    @{Row=11; Col=15; Value=3.772001065086926},
    @{Row=12; Col=2; Value=2.426853563228178},
    @{Row=12; Col=3; Value=0.2973097878531519},
    @{Row=12; Col=4; Value=0.5233465566171844},
    @{Row=12; Col=5; Value=0.1687898597026898},
    @{Row=12; Col=7; Value=0.91275346221353},
    @{Row=12; Col=8; Value=0.9515902834988736},
    @{Row=12; Col=10; Value=0.07149723185248646},
    @{Row=12; Col=12; Value=0.4829585504150344},
    @{Row=12; Col=15; Value=3.765802116040106},
    @{Row=13; Col=2; Value=2.417760366598657},
    @{Row=13; Col=3; Value=0.2963919557170698},
    @{Row=13; Col=4; Value=0.5232295430294585},
    @{Row=13; Col=5; Value=0.1688494602141066},
    @{Row=13; Col=7; Value=0.9129911036424119},
    @{Row=13; Col=8; Value=0.9519877527614824},
    @{Row=13; Col=10; Value=0.07152198493703921},
    @{Row=13; Col=12; Value=0.4821938184728509},
    @{Row=13; Col=15; Value=3.767119610049889},
    @{Row=14; Col=2; Value=2.388106635204622},
    @{Row=14; Col=3; Value=0.2933981427143522},
    @{Row=14; Col=4; Value=0.5228540703839855},
    @{Row=14; Col=5; Value=0.1690456890278664},
    @{Row=14; Col=7; Value=0.9137847329527204},
    @{Row=14; Col=8; Value=0.9532972326026226},
    @{Row=14; Col=10; Value=0.07160317561108087},
    @{Row=14; Col=12; Value=0.4797031675462335},
    @{Row=14; Col=15; Value=3.771483161770021},
    @{Row=15; Col=2; Value=2.369943799541829},
    @{Row=15; Col=3; Value=0.2915639182215841},
    @{Row=15; Col=4; Value=0.5226287972759422},
    @{Row=15; Col=5; Value=0.1691673132302967},
    @{Row=15; Col=7; Value=0.9142851629537461},
    @{Row=15; Col=8; Value=0.954109504615559},
    @{Row=15; Col=10; Value=0.07165326509021774},
    @{Row=15; Col=12; Value=0.4781801161404218},
    @{Row=15; Col=15; Value=3.774207369828105},
    @{Row=16; Col=2; Value=2.265875870509944},
    @{Row=16; Col=3; Value=0.2810463714997411},
    @{Row=16; Col=4; Value=0.5214101605502748},
    @{Row=16; Col=5; Value=0.1698861943988597},
    @{Row=16; Col=7; Value=0.9173720094053266},
    @{Row=16; Col=8; Value=0.9589202716220342},
    @{Row=16; Col=10; Value=0.07194579909449317},
    @{Row=16; Col=12; Value=0.4694912728896554},
    @{Row=16; Col=15; Value=3.790606053530524},
    @{Row=17; Col=2; Value=2.202049610652296},
    @{Row=17; Col=3; Value=0.2745887581474733},
    @{Row=17; Col=4; Value=0.5207271447165738},
    @{Row=17; Col=5; Value=0.1703467577996509},
    @{Row=17; Col=7; Value=0.9194609261275986},
    @{Row=17; Col=8; Value=0.9620106549072034},
    @{Row=17; Col=10; Value=0.07213016659250115},
    @{Row=17; Col=12; Value=0.4641960487017514},
    @{Row=17; Col=15; Value=3.80136805414304},
    @{Row=18; Col=2; Value=2.165345054775628},
    @{Row=18; Col=3; Value=0.2708726005736537},
    @{Row=18; Col=4; Value=0.5203580452797638},
    @{Row=18; Col=5; Value=0.1706188511508984},
    @{Row=18; Col=7; Value=0.9207340625498261},
    @{Row=18; Col=8; Value=0.9638392782941878},
    @{Row=18; Col=10; Value=0.07223801622168224},
    @{Row=18; Col=12; Value=0.4611633410971336},
    @{Row=18; Col=15; Value=3.807815783819933},
    @{Row=19; Col=2; Value=2.152918728915665},
    @{Row=19; Col=3; Value=0.2696140541178522},
    @{Row=19; Col=4; Value=0.5202371561699692},
    @{Row=19; Col=5; Value=0.1707122124400744},
    @{Row=19; Col=7; Value=0.9211774203872238},
    @{Row=19; Col=8; Value=0.9644671967232483},
    @{Row=19; Col=10; Value=0.0722748427652089},
    @{Row=19; Col=12; Value=0.4601387515141937},
    @{Row=19; Col=15; Value=3.810043112899479},
    @{Row=20; Col=2; Value=2.208843356143973},
    @{Row=20; Col=3; Value=0.2752763817343293},
    @{Row=20; Col=4; Value=0.5207973950228677},
    @{Row=20; Col=5; Value=0.1702969860418868},
    @{Row=20; Col=7; Value=0.9192311399002904},
    @{Row=20; Col=8; Value=0.9616763876206704},
    @{Row=20; Col=10; Value=0.07211035347496342},
    @{Row=20; Col=12; Value=0.4647583939948561},
    @{Row=20; Col=15; Value=3.800195743000586},
    @{Row=21; Col=2; Value=2.39681646488873},
    @{Row=21; Col=3; Value=0.2942775877637587},
    @{Row=21; Col=4; Value=0.5229633733720647},
    @{Row=21; Col=5; Value=0.1689877542146956},
    @{Row=21; Col=7; Value=0.9135486422119499},
    @{Row=21; Col=8; Value=0.9529104855056687},
    @{Row=21; Col=10; Value=0.07157925338277416},
    @{Row=21; Col=12; Value=0.4804342022538179},
    @{Row=21; Col=15; Value=3.770190766381603},
    @{Row=22; Col=2; Value=2.519711374222254},
    @{Row=22; Col=3; Value=0.3066770771737026},
    @{Row=22; Col=4; Value=0.5245900416631315},
    @{Row=22; Col=5; Value=0.1681960369269433},
    @{Row=22; Col=7; Value=0.9104749703953274},
    @{Row=22; Col=8; Value=0.9476369570754457},
    @{Row=22; Col=10; Value=0.0712481728853831},
    @{Row=22; Col=12; Value=0.4907932995716777},
    @{Row=22; Col=15; Value=3.752880950242314},
    @{Row=23; Col=2; Value=2.454116795484424},
    @{Row=23; Col=3; Value=0.3000610564645001},
    @{Row=23; Col=4; Value=0.5237025661143946},
    @{Row=23; Col=5; Value=0.1686127445588976},
    @{Row=23; Col=7; Value=0.9120567669398127},
    @{Row=23; Col=8; Value=0.9504098456652059},
    @{Row=23; Col=10; Value=0.07142341349889136},
    @{Row=23; Col=12; Value=0.485254086406357},
    @{Row=23; Col=15; Value=3.761908804002871},
    @{Row=24; Col=2; Value=2.205771934751908},
    @{Row=24; Col=3; Value=0.2749655182712161},
    @{Row=24; Col=4; Value=0.5207655614205891},
    @{Row=24; Col=5; Value=0.1703194650847628},
    @{Row=24; Col=7; Value=0.9193348014067482},
    @{Row=24; Col=8; Value=0.9618273481301998},
    @{Row=24; Col=10; Value=0.07211930520854715},
    @{Row=24; Col=12; Value=0.4645041215809442},
    @{Row=24; Col=15; Value=3.800724933435049},
    @{Row=25; Col=2; Value=1.938153237023187},
    @{Row=25; Col=3; Value=0.2478239885094808},
    @{Row=25; Col=4; Value=0.5185058942153091},
    @{Row=25; Col=5; Value=0.1724353567655914},
    @{Row=25; Col=7; Value=0.9299232195822924},
    @{Row=25; Col=8; Value=0.9760970936162607},
    @{Row=25; Col=10; Value=0.07293903516632305},
    @{Row=25; Col=12; Value=0.442618369632342},
    @{Row=25; Col=15; Value=3.852444053325144}
)

foreach ($item in $updates) {
    $ws.Cells.Item($item.Row, $item.Col).Value = $item.Value
}

Write-Output "Applied $($updates.Count) cell updates"